$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# Add new row 15 to the Logs sheet
$ws.Range("A15").Value = "Demo inplannen"
$ws.Range("B15").Value = "klantenservice@testbedrijf123.nl"
$ws.Range("C15").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Range("D15").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("E15").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$ws.Range("F15").Value = "2025-08-14 20:52:47"
$ws.Range("G15").Value = "Nee"
$ws.Range("H15").Value = "Ja"
$ws.Range("I15").Value = "Nee"
$ws.Range("J15").Value = "Nee"

# Extend conditional formatting ranges to include the new row (row 14 -> row 15)
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $ws.Range("$col" + "2:" + "$col" + "14")
    $newRange = $ws.Range("$col" + "2:" + "$col" + "15")
    $fc = $oldRange.FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard count for "Intern verzoek / Actie voor medewerker"
$dashboard.Range("B2").Value = 9
